# Updated cryptos list (GitHub Actions price/volume refresh).
# Numeric-looking Price (column D) values are written with a leading
# apostrophe so Excel stores them as text (matching the source data,
# which keeps Price as a string such as "68.037.50") instead of
# silently coercing them to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.037.50'
$ws.Range("E2").Value = '  +1.29%  '

$ws.Range("D3").Value = '3.894.30'
$ws.Range("E3").Value = '  -0.95%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '''486.77'
$ws.Range("E5").Value = '  +3.48%  '

$ws.Range("D6").Value = '''145.69'
$ws.Range("E6").Value = '  -0.61%  '

$ws.Range("D7").Value = '''0.618'
$ws.Range("E7").Value = '  -1.51%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = '''0.723'
$ws.Range("E9").Value = '  -1.47%  '

$ws.Range("D10").Value = '''0.164'
$ws.Range("E10").Value = '  -1.58%  '

$ws.Range("D11").Value = '''0.0000341'
$ws.Range("E11").Value = '  -0.31%  '

$ws.Range("D12").Value = '''42.57'
$ws.Range("E12").Value = '  -1.95%  '

$ws.Range("D13").Value = '''10.64'
$ws.Range("E13").Value = '  +1.62%  '

$ws.Range("D14").Value = '4.523.48'
$ws.Range("E14").Value = '  -0.70%  '

$ws.Range("D15").Value = '3.917.83'
$ws.Range("E15").Value = '  -0.17%  '

$ws.Range("B16").Value = 'TRON'
$ws.Range("C16").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D16").Value = '''0.136'
$ws.Range("E16").Value = '  -1.38%  '

$ws.Range("B17").Value = 'Uniswap'
$ws.Range("C17").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D17").Value = '''13.96'
$ws.Range("E17").Value = '  -7.99%  '

$ws.Range("D18").Value = '''19.71'
$ws.Range("E18").Value = '  -0.84%  '

$ws.Range("E19").Value = '  -3.09%  '

$ws.Range("D20").Value = '68.252.31'
$ws.Range("E20").Value = '  +1.17%  '

$ws.Range("D21").Value = '''431.41'
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").Value = '''3.51'
$ws.Range("E22").Value = '  +3.65%  '

$ws.Range("D23").Value = '''14.84'
$ws.Range("E23").Value = '  +2.42%  '

$ws.Range("D24").Value = '''87.20'
$ws.Range("E24").Value = '  -0.46%  '

$ws.Range("D25").Value = '''11.11'
$ws.Range("E25").Value = '  +13.71%  '

$ws.Range("D26").Value = '''11.29'
$ws.Range("E26").Value = '  +9.38%  '

$ws.Range("D27").Value = '''3.62'
$ws.Range("E27").Value = '  +0.76%  '

$ws.Range("D28").Value = '''37.96'
$ws.Range("E28").Value = '  -2.18%  '

$ws.Range("D29").Value = '''5.74'
$ws.Range("E29").Value = '  -0.22%  '

$ws.Range("D30").Value = '''721.21'
$ws.Range("E30").Value = '  +0.00%  '

$ws.Range("D31").Value = '''13.59'
$ws.Range("E31").Value = '  +0.42%  '

$ws.Range("D32").Value = '''0.129'
$ws.Range("E32").Value = '  -2.64%  '

$ws.Range("D33").Value = '''2.89'
$ws.Range("E33").Value = '  +2.60%  '

$ws.Range("D34").Value = '''6.23'
$ws.Range("E34").Value = '  +16.15%  '

$ws.Range("D35").Value = '''41.35'
$ws.Range("E35").Value = '  -3.00%  '

$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").Value = '0.0₃0856'
$ws.Range("E36").Value = '  +6.51%  '

$ws.Range("B37").Value = 'OKB'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D37").Value = '''59.90'
$ws.Range("E37").Value = '  +3.42%  '

$ws.Range("D38").Value = '''0.404'
$ws.Range("E38").Value = '  +19.64%  '

$ws.Range("D39").Value = '''0.147'
$ws.Range("E39").Value = '  -2.69%  '

$ws.Range("D40").Value = '''0.999'
$ws.Range("E40").Value = '  +0.02%  '

$ws.Range("D41").Value = '''2.95'
$ws.Range("E41").Value = '  +15.05%  '

$ws.Range("D42").Value = '''0.0476'
$ws.Range("E42").Value = '  -0.33%  '

$ws.Range("D43").Value = '''3.12'
$ws.Range("E43").Value = '  +2.17%  '

$ws.Range("D44").Value = '''2.90'
$ws.Range("E44").Value = '  +2.89%  '

$ws.Range("B45").Value = 'FirstDigitalUSD'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D45").Value = '''1.00'
$ws.Range("E45").Value = '  +0.12%  '

$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").Value = '''0.140'
$ws.Range("E46").Value = '  -1.30%  '

$ws.Range("D47").Value = '''3.28'
$ws.Range("E47").Value = '  +2.84%  '

$ws.Range("D48").Value = '''3.37'
$ws.Range("E48").Value = '  -5.07%  '

$ws.Range("D49").Value = '''2.11'
$ws.Range("E49").Value = '  -4.08%  '

$ws.Range("D50").Value = '''144.79'
$ws.Range("E50").Value = '  -1.91%  '

$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").Value = '''2.80'
$ws.Range("E51").Value = '  -2.67%  '
